# Applies the "43 特殊符号查询表" edit:
#  - P1: merge the title's two runs into one (drop the mid-paragraph
#    "_GoBack" bookmark) and append "python3的特殊符号查询表" in place.
#  - Re-flow / comment several source-listing paragraphs.
#  - Paragraph 8 gains a second run (Tab + explanation text) and the
#    "_GoBack" bookmark is re-created at the end of that paragraph.
#  - One previously-empty paragraph gains a "#定义参数处理类型" comment run.
#  - A handful of one-line text tweaks (trailing "#comment" appended).

$d = $word.ActiveDocument

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Set-ParaXml($paraIndex, $bodyXml) {
    $rng = $d.Paragraphs($paraIndex).Range
    $pkg = "<pkg:package xmlns:pkg=`"http://schemas.microsoft.com/office/2006/xmlPackage`">" +
           "<pkg:part pkg:name=`"/word/document.xml`" pkg:contentType=`"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml`">" +
           "<pkg:xmlData><w:document $wNs><w:body>$bodyXml</w:body></w:document></pkg:xmlData>" +
           "</pkg:part></pkg:package>"
    $rng.InsertXML($pkg)
}

# --- Paragraph 1: "根据...构造一个基于" + bookmark + "python3的特殊符号查询表" ---
# becomes a single run with the bookmark removed.
$p1Body = '<w:p><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr></w:pPr>' +
          '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr>' +
          '<w:t>根据2017年1月8日的一项计划，构造一个基于python3的特殊符号查询表</w:t></w:r></w:p>'
Set-ParaXml 1 $p1Body

# --- Paragraph 5: "help='''" -> "help=''' #定义帮助信息" ---
$d.Paragraphs(5).Range.Text = "help=''' #定义帮助信息"

# --- Paragraph 8: "    findspec.py --no-original|-j theta" gains a second
#     run (Tab + explanation) and the "_GoBack" bookmark moves here. ---
$p8Body = '<w:p><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr></w:pPr>' +
          '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr>' +
          '<w:t xml:space="preserve">    findspec.py --no-original|-j theta</w:t></w:r>' +
          '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr>' +
          '<w:tab/><w:t>, means do not print original name of the symbol</w:t></w:r>' +
          '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'
Set-ParaXml 8 $p8Body

# --- Paragraph 24: previously-empty paragraph gains a comment run ---
$p24Body = '<w:p><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr></w:pPr>' +
           '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr>' +
           '<w:t>#定义参数处理类型</w:t></w:r></w:p>'
Set-ParaXml 24 $p24Body

# --- Paragraph 35: "args = sys.argv[1:]" -> append comment ---
$d.Paragraphs(35).Range.Text = "args = sys.argv[1:] #获取命令行参数"

# --- Paragraph 43: "values = {" -> append comment ---
$d.Paragraphs(43).Range.Text = "values = { #获取解析参数"

# --- Paragraph 50: short_options -> append comment (keep straight quotes) ---
$d.Paragraphs(50).Range.Text = 'short_options = "n:sjha" #定义长短选项'

# --- Paragraph 62: if values["do_help"]!=None: -> append comment ---
$d.Paragraphs(62).Range.Text = 'if values["do_help"]!=None: #开始处理过程，一般先处理help选项'
